$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 166.71428
$ws.Range("I5").Value = 35
$ws.Range("J5").Value = 342.33334
$ws.Range("K5").Value = 35
$ws.Range("L5").Value = 342.33334
$ws.Range("M5").Value = 80
$ws.Range("N5").Value = -572.33334
$ws.Range("H41").Value = 258
$ws.Range("I41").Value = 176
$ws.Range("K41").Value = 176
$ws.Range("M41").Value = 264
$ws.Range("H55").Value = 278.6
$ws.Range("I55").Value = 381
$ws.Range("K55").Value = 381
$ws.Range("M55").Value = -167
$ws.Range("H62").Value = 8230.667
$ws.Range("I62").Value = 8237
$ws.Range("K62").Value = 8237
$ws.Range("M62").Value = -7613
$ws.Range("H65").Value = 8230.667
$ws.Range("I65").Value = 8237
$ws.Range("K65").Value = 41185
$ws.Range("M65").Value = -38065
$ws.Range("H74").Value = 7073
$ws.Range("I74").Value = 7073
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 7073
$ws.Range("L74").Value = 0
$ws.Range("M74").ClearContents()
$ws.Range("N74").Value = -6137
$ws.Range("H77").Value = 7073
$ws.Range("I77").Value = 7073
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 35365
$ws.Range("L77").Value = 0
$ws.Range("M77").ClearContents()
$ws.Range("N77").Value = -30685
$ws.Range("H96").Value = 335.66666
$ws.Range("I96").Value = 218.15384
$ws.Range("J96").Value = 1099.5
$ws.Range("K96").Value = 654.4615200000001
$ws.Range("L96").Value = 3298.5
$ws.Range("M96").Value = 718.5384799999999
$ws.Range("N96").Value = -6044.5
$ws.Range("H116").Value = 4474.6665
$ws.Range("I116").Value = 4466.5
$ws.Range("J116").Value = 4491
$ws.Range("K116").Value = 4466.5
$ws.Range("L116").Value = 4491
$ws.Range("M116").Value = -1024.5
$ws.Range("N116").Value = -11375

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3439.8
$ws.Range("I2").Value = 3439.8
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 3439.8
$ws.Range("L2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("N2").Value = -3326.8
$ws.Range("H32").Value = 10309.315
$ws.Range("I32").Value = 10826.5
$ws.Range("K32").Value = 10826.5
$ws.Range("M32").Value = -10539.5
$ws.Range("H116").Value = 3439.8
$ws.Range("I116").Value = 3439.8
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 3439.8
$ws.Range("L116").Value = 0
$ws.Range("M116").ClearContents()
$ws.Range("N116").Value = -1145.8

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3439.8
$ws.Range("I3").Value = 3439.8
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 3439.8
$ws.Range("L3").Value = 0
$ws.Range("M3").ClearContents()
$ws.Range("N3").Value = -3325.8

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 1383.6842
$ws.Range("I107").Value = 1545.3334
$ws.Range("J107").Value = 777.5
$ws.Range("K107").Value = 1545.3334
$ws.Range("L107").Value = 777.5
$ws.Range("M107").Value = 374.6666
$ws.Range("N107").Value = -4617.5
$ws.Range("H122").Value = 5166.6665
$ws.Range("I122").Value = 7200.3335
$ws.Range("J122").Value = 3133
$ws.Range("K122").Value = 21601.0005
$ws.Range("L122").Value = 9399
$ws.Range("M122").Value = -19151.0005
$ws.Range("N122").Value = -14299

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 1717554.9
$ws.Range("I8").Value = 1717554.9
$ws.Range("K8").Value = 5152664.699999999
$ws.Range("M8").Value = -5152525.699999999
$ws.Range("H22").Value = 1000
$ws.Range("I22").Value = 1000
$ws.Range("K22").Value = 3000
$ws.Range("M22").Value = -2831
$ws.Range("H27").Value = 1000
$ws.Range("I27").Value = 1000
$ws.Range("K27").Value = 3000
$ws.Range("M27").Value = -2898
$ws.Range("H40").Value = 353.2
$ws.Range("I40").Value = 322
$ws.Range("K40").Value = 1288
$ws.Range("M40").Value = -1219
$ws.Range("H86").Value = 1537.5714
$ws.Range("I86").Value = 1818
$ws.Range("J86").Value = 1425.4
$ws.Range("K86").Value = 5454
$ws.Range("L86").Value = 4276.200000000001
$ws.Range("M86").Value = -4268
$ws.Range("N86").Value = -6648.200000000001
$ws.Range("H89").Value = 1537.5714
$ws.Range("I89").Value = 1818
$ws.Range("J89").Value = 1425.4
$ws.Range("K89").Value = 16362
$ws.Range("L89").Value = 12828.6
$ws.Range("M89").Value = -10434
$ws.Range("N89").Value = -24684.6
$ws.Range("H92").Value = 548.1667
$ws.Range("J92").Value = 694.5
$ws.Range("L92").Value = 2083.5
$ws.Range("N92").Value = -4579.5
$ws.Range("H113").Value = 1234
$ws.Range("I113").Value = 595.75
$ws.Range("K113").Value = 1787.25
$ws.Range("M113").Value = 382.75
$ws.Range("H122").Value = 3072.7856
$ws.Range("I122").Value = 1041.6666
$ws.Range("J122").Value = 3626.7273
$ws.Range("K122").Value = 9374.9994
$ws.Range("L122").Value = 32640.5457
$ws.Range("M122").Value = -6924.999400000001
$ws.Range("N122").Value = -37540.5457
$ws.Range("H132").Value = 1398.8
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()
$ws.Range("H133").Value = 0
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("L133").ClearContents()
$ws.Range("M133").ClearContents()
$ws.Range("N133").Value = 0

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 80.875
$ws.Range("I2").Value = 77.333336
$ws.Range("J2").Value = 91.5
$ws.Range("K2").Value = 77.333336
$ws.Range("L2").Value = 91.5
$ws.Range("M2").Value = 35.666664
$ws.Range("N2").Value = -317.5
$ws.Range("H107").Value = 4666.3335
$ws.Range("I107").Value = 4666.3335
$ws.Range("K107").Value = 4666.3335
$ws.Range("M107").Value = -2746.3335
$ws.Range("H126").Value = 4400
$ws.Range("I126").Value = 4750
$ws.Range("K126").Value = 14250
$ws.Range("M126").Value = -11780

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 4867.364
$ws.Range("I22").Value = 3243
$ws.Range("J22").Value = 7710
$ws.Range("K22").Value = 3243
$ws.Range("L22").Value = 7710
$ws.Range("M22").Value = -2948
$ws.Range("N22").Value = -8300
$ws.Range("H27").Value = 4867.364
$ws.Range("I27").Value = 3243
$ws.Range("J27").Value = 7710
$ws.Range("K27").Value = 3243
$ws.Range("L27").Value = 7710
$ws.Range("M27").Value = -3136
$ws.Range("N27").Value = -7924
$ws.Range("H40").Value = 3247
$ws.Range("I40").Value = 1494
$ws.Range("K40").Value = 1494
$ws.Range("M40").Value = -1358

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 40000
$ws.Range("J32").Value = 40000
$ws.Range("L32").Value = 40000
$ws.Range("N32").Value = -40634
$ws.Range("H54").Value = 32634.857
$ws.Range("J54").Value = 32634.857
$ws.Range("L54").Value = 32634.857
$ws.Range("N54").Value = -33674.857
$ws.Range("H81").Value = 2635.2856
$ws.Range("I81").Value = 2635.2856
$ws.Range("K81").Value = 5270.5712
$ws.Range("M81").Value = -4209.5712
$ws.Range("H84").Value = 2635.2856
$ws.Range("I84").Value = 2635.2856
$ws.Range("K84").Value = 26352.856
$ws.Range("M84").Value = -21048.856
$ws.Range("H113").Value = 464.375
$ws.Range("I113").Value = 525.8333
$ws.Range("J113").Value = 280
$ws.Range("K113").Value = 1577.4999
$ws.Range("L113").Value = 840
$ws.Range("M113").Value = 592.5001
$ws.Range("N113").Value = -5180
$ws.Range("H122").Value = 1641.25
$ws.Range("J122").Value = 1988
$ws.Range("L122").Value = 5964
$ws.Range("N122").Value = -10864
$ws.Range("H126").Value = 1885.7142
